# Generate Report for Archive
#
# 1. The status text "Ready for handoff" changes to "In Translation"
#    (shows up in the Overview sheet's zh-cn/de-de status cells, and in
#    the Status column of the zh-cn / de-de detail sheets).
# 2. Because the new status text is shorter, the corresponding status
#    columns are narrowed from ~17.22 width to ~13.41 width.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Narrow the status columns to match the new shorter text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
